# Natmi following Dr Hou advice
#
# The sending/target cluster set grows from {M2} to {M2, sCs}, so the
# ligand-receptor pair table (Ly9 -> Ly9) now has one row per
# (sending cluster, target cluster) combination: 2 x 2 = 4 rows instead
# of the original single M2->M2 row. Row 2 is updated in place and three
# new rows (3-5) are appended below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; A = "M2";  B = "Ly9"; C = "Ly9"; D = "M2";
       E = 3; F = 1;                  G = 44.76975266666667; H = 134.309258;
       I = 0.9848052097546599;  J = 0.9848052097546598;
       K = 3; L = 1;                  M = 44.76975266666667; N = 134.309258;
       O = 0.9848052097546599;  P = 0.9848052097546598;
       Q = 2004.330753834507;   R = 18038.97678451056;
       S = 0.9698413011599197;  T = 0.9698413011599195 },

    @{ Row = 3; A = "M2";  B = "Ly9"; C = "Ly9"; D = "sCs";
       E = 3; F = 1;                  G = 44.76975266666667; H = 134.309258;
       I = 0.9848052097546599;  J = 0.9848052097546598;
       K = 1; L = 0.3333333333333333; M = 0.690763;          N = 2.072289;
       O = 0.01519479024534016; P = 0.01519479024534015;
       Q = 30.92528866128467;   R = 278.327597951562;
       S = 0.01496390859474027; T = 0.01496390859474027 },

    @{ Row = 4; A = "sCs"; B = "Ly9"; C = "Ly9"; D = "M2";
       E = 1; F = 0.3333333333333333; G = 0.690763;          H = 2.072289;
       I = 0.01519479024534016; J = 0.01519479024534015;
       K = 3; L = 1;                  M = 44.76975266666667; N = 134.309258;
       O = 0.9848052097546599;  P = 0.9848052097546598;
       Q = 30.92528866128467;   R = 278.327597951562;
       S = 0.01496390859474027; T = 0.01496390859474027 },

    @{ Row = 5; A = "sCs"; B = "Ly9"; C = "Ly9"; D = "sCs";
       E = 1; F = 0.3333333333333333; G = 0.690763;          H = 2.072289;
       I = 0.01519479024534016; J = 0.01519479024534015;
       K = 1; L = 0.3333333333333333; M = 0.690763;          N = 2.072289;
       O = 0.01519479024534016; P = 0.01519479024534015;
       Q = 0.477153522169;      R = 4.294381699521;
       S = 0.0002308816505998844; T = 0.0002308816505998843 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($rowData in $data) {
    $r = $rowData.Row
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}
